$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(79).Insert()

$ws.Range("A79").Value = 3
$ws.Range("B79").Value = "Femacal de La Calera"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44719
$ws.Range("E79").Value = 5
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100107
$ws.Range("H79").Value = "Otros"
$ws.Range("I79").Value = 100107011
$ws.Range("J79").Value = "Tuna"
$ws.Range("K79").Value = "Sin especificar"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 67
$ws.Range("N79").Value = 17000
$ws.Range("O79").Value = 17000
$ws.Range("P79").Value = 17000
$ws.Range("Q79").Value = '$/caja 20 kilos'
$ws.Range("R79").Value = "Provincia de Limarí"
$ws.Range("S79").Value = 850
$ws.Range("T79").Value = 20
